$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPJKGHBXD%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJHMEUCIEL3AfXLHDjJrmiZQ7nvts%2F2%2FexeP1deyhI%2B869CilBCAiEA4a35CulIz8PILqKQLVTaZjMw472mi%2Fp4qMEiw9hOAk4q3QMIGBAEGgw3ODEyMDQ5NDIyNDQiDH0vN0EMjcsYmXsWXCq6A7F2ywjulPzTT%2Fl65V8dHrEqkMzLSLzsFj3IU3M2HXCLQ0nY%2BUxKRgzAtSVe40O8cinK%2B%2FEc7yOTvjyvox%2FUGZUB%2BoBxrbtmgKTWoZRKqHWnoh%2FMGM9EG71QmeOtMVn0FdD%2Fq78hCsmdvGaSkwGSb4RbtVts9JRU51Uo7nm7apnf%2BIWawVSrdqZ37PZV7VItmAyPZbjaFiTNOt0p6vQ5rm7qjSCoK4Xc%2Fsn2Q%2BE4cXSny9IajGE1Tglvd9Al1SZzptwu9jDtfWleVfKXeh4MNKRlLHli7icz4m7Zktu1ac8J8NslNWhHurpWoUeI967GldsanCMQq0e0%2BEaRsT7TfNqyWqQlPo%2B5M%2FFqTBdU1LexNr4D%2B7t%2FAJsPPHTu7jDiNdQL0H10UQLtV2YgMR3YGXA9Crq4Rhc5P90xuJmZuLxQhRZquzktqt6XAESwcax2vnq6%2Fg9Kc%2F1PUGqHprHUde90Gj%2BluPafSB7orTUXIo0RviWahzfSXEvwlHZNUcTWSOu6HeTTWdRCjkQxs1C5u9s4nKhJbuRjryFawVDAJyuICs3KYs5uCWA0KCCqnWF1zNZ9oZAT90aLbKowyZX8wQY6pQFC5qd2nN425EfejESqT6OdYkJSoxJWIE%2FPFuQ4XbZdPPGTZPSUhuxnzuzAb%2FKPMs31zAGUd7akWHflaehpZ6nzWr17U%2Bd0v3B3mPEn6Lt3fqbtyZGEg3TElvrX4Qt%2B7C94LXQVYw2NJSBLTE2dR9%2B%2BrPhvIvsEHdKciWuXjSRUgfQKAZTKyKelZ9F%2FRbXKFxlhCaB7N9r5jUAKFPwPA67C9ZRgwgI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=e1083626b499b2f80f41f5a31a7c14b667bdf16b634edaf33b4738be54d10234'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPJKGHBXD%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJHMEUCIEL3AfXLHDjJrmiZQ7nvts%2F2%2FexeP1deyhI%2B869CilBCAiEA4a35CulIz8PILqKQLVTaZjMw472mi%2Fp4qMEiw9hOAk4q3QMIGBAEGgw3ODEyMDQ5NDIyNDQiDH0vN0EMjcsYmXsWXCq6A7F2ywjulPzTT%2Fl65V8dHrEqkMzLSLzsFj3IU3M2HXCLQ0nY%2BUxKRgzAtSVe40O8cinK%2B%2FEc7yOTvjyvox%2FUGZUB%2BoBxrbtmgKTWoZRKqHWnoh%2FMGM9EG71QmeOtMVn0FdD%2Fq78hCsmdvGaSkwGSb4RbtVts9JRU51Uo7nm7apnf%2BIWawVSrdqZ37PZV7VItmAyPZbjaFiTNOt0p6vQ5rm7qjSCoK4Xc%2Fsn2Q%2BE4cXSny9IajGE1Tglvd9Al1SZzptwu9jDtfWleVfKXeh4MNKRlLHli7icz4m7Zktu1ac8J8NslNWhHurpWoUeI967GldsanCMQq0e0%2BEaRsT7TfNqyWqQlPo%2B5M%2FFqTBdU1LexNr4D%2B7t%2FAJsPPHTu7jDiNdQL0H10UQLtV2YgMR3YGXA9Crq4Rhc5P90xuJmZuLxQhRZquzktqt6XAESwcax2vnq6%2Fg9Kc%2F1PUGqHprHUde90Gj%2BluPafSB7orTUXIo0RviWahzfSXEvwlHZNUcTWSOu6HeTTWdRCjkQxs1C5u9s4nKhJbuRjryFawVDAJyuICs3KYs5uCWA0KCCqnWF1zNZ9oZAT90aLbKowyZX8wQY6pQFC5qd2nN425EfejESqT6OdYkJSoxJWIE%2FPFuQ4XbZdPPGTZPSUhuxnzuzAb%2FKPMs31zAGUd7akWHflaehpZ6nzWr17U%2Bd0v3B3mPEn6Lt3fqbtyZGEg3TElvrX4Qt%2B7C94LXQVYw2NJSBLTE2dR9%2B%2BrPhvIvsEHdKciWuXjSRUgfQKAZTKyKelZ9F%2FRbXKFxlhCaB7N9r5jUAKFPwPA67C9ZRgwgI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=f6ec20c9599a59c34c99cebcb8c75ba2ad80ad76ebd2c699f53d78a18f036d66'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPJKGHBXD%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJHMEUCIEL3AfXLHDjJrmiZQ7nvts%2F2%2FexeP1deyhI%2B869CilBCAiEA4a35CulIz8PILqKQLVTaZjMw472mi%2Fp4qMEiw9hOAk4q3QMIGBAEGgw3ODEyMDQ5NDIyNDQiDH0vN0EMjcsYmXsWXCq6A7F2ywjulPzTT%2Fl65V8dHrEqkMzLSLzsFj3IU3M2HXCLQ0nY%2BUxKRgzAtSVe40O8cinK%2B%2FEc7yOTvjyvox%2FUGZUB%2BoBxrbtmgKTWoZRKqHWnoh%2FMGM9EG71QmeOtMVn0FdD%2Fq78hCsmdvGaSkwGSb4RbtVts9JRU51Uo7nm7apnf%2BIWawVSrdqZ37PZV7VItmAyPZbjaFiTNOt0p6vQ5rm7qjSCoK4Xc%2Fsn2Q%2BE4cXSny9IajGE1Tglvd9Al1SZzptwu9jDtfWleVfKXeh4MNKRlLHli7icz4m7Zktu1ac8J8NslNWhHurpWoUeI967GldsanCMQq0e0%2BEaRsT7TfNqyWqQlPo%2B5M%2FFqTBdU1LexNr4D%2B7t%2FAJsPPHTu7jDiNdQL0H10UQLtV2YgMR3YGXA9Crq4Rhc5P90xuJmZuLxQhRZquzktqt6XAESwcax2vnq6%2Fg9Kc%2F1PUGqHprHUde90Gj%2BluPafSB7orTUXIo0RviWahzfSXEvwlHZNUcTWSOu6HeTTWdRCjkQxs1C5u9s4nKhJbuRjryFawVDAJyuICs3KYs5uCWA0KCCqnWF1zNZ9oZAT90aLbKowyZX8wQY6pQFC5qd2nN425EfejESqT6OdYkJSoxJWIE%2FPFuQ4XbZdPPGTZPSUhuxnzuzAb%2FKPMs31zAGUd7akWHflaehpZ6nzWr17U%2Bd0v3B3mPEn6Lt3fqbtyZGEg3TElvrX4Qt%2B7C94LXQVYw2NJSBLTE2dR9%2B%2BrPhvIvsEHdKciWuXjSRUgfQKAZTKyKelZ9F%2FRbXKFxlhCaB7N9r5jUAKFPwPA67C9ZRgwgI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=f123ad0a4d34c4ba76b08f0e509b336dbb6590ec261ebd87d028adaf7f8272ec'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPJKGHBXD%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T180033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJHMEUCIEL3AfXLHDjJrmiZQ7nvts%2F2%2FexeP1deyhI%2B869CilBCAiEA4a35CulIz8PILqKQLVTaZjMw472mi%2Fp4qMEiw9hOAk4q3QMIGBAEGgw3ODEyMDQ5NDIyNDQiDH0vN0EMjcsYmXsWXCq6A7F2ywjulPzTT%2Fl65V8dHrEqkMzLSLzsFj3IU3M2HXCLQ0nY%2BUxKRgzAtSVe40O8cinK%2B%2FEc7yOTvjyvox%2FUGZUB%2BoBxrbtmgKTWoZRKqHWnoh%2FMGM9EG71QmeOtMVn0FdD%2Fq78hCsmdvGaSkwGSb4RbtVts9JRU51Uo7nm7apnf%2BIWawVSrdqZ37PZV7VItmAyPZbjaFiTNOt0p6vQ5rm7qjSCoK4Xc%2Fsn2Q%2BE4cXSny9IajGE1Tglvd9Al1SZzptwu9jDtfWleVfKXeh4MNKRlLHli7icz4m7Zktu1ac8J8NslNWhHurpWoUeI967GldsanCMQq0e0%2BEaRsT7TfNqyWqQlPo%2B5M%2FFqTBdU1LexNr4D%2B7t%2FAJsPPHTu7jDiNdQL0H10UQLtV2YgMR3YGXA9Crq4Rhc5P90xuJmZuLxQhRZquzktqt6XAESwcax2vnq6%2Fg9Kc%2F1PUGqHprHUde90Gj%2BluPafSB7orTUXIo0RviWahzfSXEvwlHZNUcTWSOu6HeTTWdRCjkQxs1C5u9s4nKhJbuRjryFawVDAJyuICs3KYs5uCWA0KCCqnWF1zNZ9oZAT90aLbKowyZX8wQY6pQFC5qd2nN425EfejESqT6OdYkJSoxJWIE%2FPFuQ4XbZdPPGTZPSUhuxnzuzAb%2FKPMs31zAGUd7akWHflaehpZ6nzWr17U%2Bd0v3B3mPEn6Lt3fqbtyZGEg3TElvrX4Qt%2B7C94LXQVYw2NJSBLTE2dR9%2B%2BrPhvIvsEHdKciWuXjSRUgfQKAZTKyKelZ9F%2FRbXKFxlhCaB7N9r5jUAKFPwPA67C9ZRgwgI%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2002%3A00%3A33%20GMT&X-Amz-Signature=17681991d21d68de45c5d5ea4aad14e57741e567ea063eee1cd9b4b96c215e33'

$wb.Save()
